$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Part 1" to "Part 2"
$ws.Name = "Part 2"

# Start from a clean sheet and rewrite the (now shorter) layout so the row
# numbers/shared-string table line up exactly with the target shape instead
# of relying on fragile row insert/delete shifting.
$ws.Cells.Clear()

# Inputs section labels (column A)
$ws.Range("A1").Value = "Inputs"
$ws.Range("A2").Value = "Arrival Rate"
$ws.Range("A3").Value = "Service Rate"
$ws.Range("A4").Value = "P(W > 0) Less Than"
$ws.Range("A5").Value = "E(W) Less Than"

# Inputs section values (column B) -- these are numeric-looking text, so
# force the "Text" number format before writing them so they stay strings.
$ws.Range("B1:B5").NumberFormat = "@"
$ws.Range("B1").Value = "Values"
$ws.Range("B2").Value = "10.0"
$ws.Range("B3").Value = "11.0"
$ws.Range("B4").Value = "0.2"
$ws.Range("B5").Value = "0.25"

# Results section labels
$ws.Range("A6").Value = "Results"
$ws.Range("A7").Value = "Number of Servers"

# Results section value
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "2"
